# Re-create documentation after refactoring: append the "Task 15 - Summing up"
# section (and relocate the _GoBack bookmark into its own trailing paragraph)
# at the end of the document.

$d = $word.ActiveDocument

# WordprocessingML namespace used for the raw paragraph fragments below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Small helper: build a <w:p> fragment for a single-run paragraph.
function New-SimpleParaXml([string]$text, [bool]$bold) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    if ($bold) {
        return "<w:p $wNs><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>$escaped</w:t></w:r></w:p>"
    }
    return "<w:p $wNs><w:r><w:t>$escaped</w:t></w:r></w:p>"
}

# The _GoBack bookmark currently sits at the very end of the "Task 11" paragraph.
# Remove it from there; it will be re-inserted into its own empty paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Build the ordered list of new paragraphs to append.
$paragraphs = @()

$paragraphs += New-SimpleParaXml "Task 15 - Summing up" $true
$paragraphs += New-SimpleParaXml "What I learned" $true

$paragraphs += New-SimpleParaXml ("As projects get bigger, the amount of code smells in them tends to increase. " + `
    "The Monopoly code had a lot of bad smells, while the code for project 1 only had 1 bad smell. " + `
    "This is expected as larger code bases are often developed with a larger team and not enough time is allotted for refactoring. " + `
    "However, refactoring a large code base after a long period of development often results in too many code smells. " + `
    "Refactoring a large number of code smells can be tough. It might be a better idea to keep looking for code smells as the code is developed.") $false

# This paragraph has two runs split by a lastRenderedPageBreak, as in the source document.
$manyPart1 = ("Many of the code smells directly affect future development. For example, fixing many of the " + `
    [char]0x201C + "feature envy" + [char]0x201D + " code smells often result in improved classes with higher cohesion. " + `
    "This allows for better reuse of existing classes and reduces the amount of duplicate code moving forward. " + `
    "However, it also shows that ") -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
$manyPart2 = ("maintenance of code should be an ongoing effort. Refactoring after a long period of development can often get too cumbersome and in turn make code harder to maintain.") -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
$paragraphs += "<w:p $wNs><w:r><w:t xml:space=`"preserve`">$manyPart1</w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>$manyPart2</w:t></w:r></w:p>"

$paragraphs += New-SimpleParaXml "Eclipse/JDeodorant" $true

# This paragraph also has two runs (split mid-word, as in the source document).
$eclipsePart1 = ("Eclipse has great refactoring tools which help a developer improve their code without too much effort. " + `
    "Most of the refactoring often takes care of any possible problems which can be caused by refactoring. " + `
    "For example, manually changing the name of a method in where it is declared would result in errors in classes where the method is used. " + `
    "Refactoring automatically changes the name of the method in all places where it is cal") -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
$eclipsePart2 = ("led. Functionality like this allows the developer to constantly refactor their code and make it more maintainable in the future.") -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
$paragraphs += "<w:p $wNs><w:r><w:t>$eclipsePart1</w:t></w:r><w:r><w:t>$eclipsePart2</w:t></w:r></w:p>"

# Empty paragraph that now owns the relocated _GoBack bookmark.
$paragraphs += "<w:p $wNs><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$paragraphs += New-SimpleParaXml "Unit tests" $true

$fragmentXml = $paragraphs -join ''

# Insert everything in one shot at the (freshly resolved) end of the document,
# after the existing "Task 11" paragraph.
$insertionPoint = $d.Content.End
$target = $d.Range($insertionPoint, $insertionPoint)
[void]$target.InsertXML($fragmentXml)

Write-Output "Appended Task 15 section ($($paragraphs.Count) paragraphs)."
